## reg_RMSE.xlsx - update regressor RMSE coefficients on the "UK" sheet
## (commit: "incomplete validation of new social care routines")
##
## S2k (row 24) RMSE:  0.94706000000000001 -> 0.94330000000000003
## S3e (row 25) RMSE:  1.2428999999999999  -> 1.2788999999999999
##
## Also replays the author's final view/selection state on the UK sheet
## (scrolled down so row 14 is the top visible row, with cell O25 selected).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("UK")
$ws.Activate()

# --- updated coefficients -------------------------------------------------
$ws.Range("B24").Value = 0.94330000000000003
$ws.Range("B25").Value = 1.2788999999999999

# --- view state: scroll so row 14 becomes the top-left visible row -------
$excel.ActiveWindow.ScrollRow = 14
$excel.ActiveWindow.ScrollColumn = 1

# --- final selection -------------------------------------------------------
$ws.Range("O25").Select()
